# "working on date issues"
# A new reporter/partner record ("Canada") was missing from the data and
# needs to be inserted as row 23 (between "Cambodia" and the existing
# "Cayman Isds" row), pushing every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 23. Excel shifts rows 23:142 down to
# 24:143 (carrying their values/formatting with them) and extends the used
# range accordingly.
$ws.Rows("23").Insert()

# Fill in the newly inserted row with the missing record.
$ws.Range("A23").Value = 44197
$ws.Range("B23").Value = 501538854874.334
$ws.Range("C23").Value = "Canada"
$ws.Range("D23").Value = "World"
